# Apply the "Media and Wiring" connector short-name relabeling and add
# chamber-related connector rows to the Connectors table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename several "Short Name" values in column A (rows 2-11) ---
$ws.Range("A2").Value  = "HE"    # was H0 (Hotend)
$ws.Range("A3").Value  = "TH"    # was TH0 (Hotend Thermistor)
$ws.Range("A4").Value  = "FH1"   # was F0 (CNC Fan/Hotend Fan)
$ws.Range("A5").Value  = "FH2"   # was F1 (Part Cooling Fan)
$ws.Range("A7").Value  = "ME"    # was E0 (Extruder Stepper)
$ws.Range("A8").Value  = "SX"    # was X (X Limit Switch)
$ws.Range("A9").Value  = "SY"    # was Y (Y Limit Switch)
$ws.Range("A10").Value = "SZ1"   # was Z1 (Left Z Limit Switch)
$ws.Range("A11").Value = "SZ2"   # was Z2 (Right Z Limit Switch)

# --- Bed Power / Bed Thermistor rows ---
$ws.Range("A15").Value = "BED"   # was HB (Bed Power)
$ws.Range("A16").Value = "TB"    # was THB (Bed Thermistor)
$ws.Range("C16").Value = "TB"    # was THB

# --- Insert three new rows for chamber connectors, just above
#     "External Connectors" (currently row 17) ---
$ws.Rows("17:19").Insert()

$ws.Range("A17").Value = "TC"
$ws.Range("B17").Value = "Chamber Thermistor"

$ws.Range("A18").Value = "FC1"
$ws.Range("B18").Value = "Chamber Exhaust"

$ws.Range("A19").Value = "FC2"
$ws.Range("B19").Value = "Chamber Filtration"

# --- Rename the "External Connectors" sub-table short names
#     (now at rows 21-24 after the insert) ---
$ws.Range("A21").Value = "MA"    # was MOTA (Alpha Stepper)
$ws.Range("A22").Value = "MB"    # was MOTB (Beta Stepper)
$ws.Range("A23").Value = "MZ1"   # was MOTZ1 (Z1 Stepper)
$ws.Range("A24").Value = "MZ2"   # was MOTZ2 (Z2 Stepper)

# --- Resize the table / autofilter to cover the three new rows ---
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:C24"))

# --- Restore the active selection to C13 ---
$ws.Range("C13").Select()
